$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.215.88"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "2.180.31"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("E6").Value = "  -3.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.71"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.22%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.14"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.84"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -9.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0936"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.04%  "

$ws.Range("D15").Value = "2.505.34"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.48%  "

$ws.Range("D18").Value = "2.145.80"
$ws.Range("E18").Value = "  -3.70%  "

$ws.Range("D19").Value = "41.141.17"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  -1.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.68"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.25%  "

$ws.Range("E25").Value = "  -5.41%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.94%  "

$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.97"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.61%  "

$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0750"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.21%  "

$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.52"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.62"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.93%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.58"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +15.47%  "

$ws.Range("E41").Value = "  -3.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.74%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.52%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.189"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.12%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0992"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.01%  "

$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  -2.76%  "

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.28%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.90%  "
